$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("R1").Value = "Metered Status"
$ws.Range("A1:R1").NumberFormat = "@"
